# ---------------------------------------------------------------
# db_check.xlsx edit: add sheet2/sheet21/sheet22 (rejected authors)
# plus Date_Created / Date_Expired / extra URL column on Sheet1.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Update Sheet1: populate new C (Date_Created) / D (Date_Expired) / G (URL1 decision) columns ---
$ws1.Cells.Item(2,3).Value = 43760.61611694683
$ws1.Cells.Item(2,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(2,4).Value = 44126.61611694683
$ws1.Cells.Item(2,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(2,7).Value = "No, individual is not listed"
$ws1.Cells.Item(3,3).Value = 43760.61614441942
$ws1.Cells.Item(3,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(3,4).Value = 44126.61614441942
$ws1.Cells.Item(3,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(3,7).Value = "No, individual is not listed"
$ws1.Cells.Item(4,3).Value = 43760.61617156792
$ws1.Cells.Item(4,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(4,4).Value = 44126.61617156792
$ws1.Cells.Item(4,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(4,7).Value = "No, individual is not listed"
$ws1.Cells.Item(5,3).Value = 43760.61619973052
$ws1.Cells.Item(5,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(5,4).Value = 44126.61619973052
$ws1.Cells.Item(5,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(5,7).Value = "No, individual is not listed"
$ws1.Cells.Item(6,3).Value = 43760.6162252288
$ws1.Cells.Item(6,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(6,4).Value = 44126.6162252288
$ws1.Cells.Item(6,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(6,7).Value = "No, individual is not listed"
$ws1.Cells.Item(7,3).Value = 43760.61624953821
$ws1.Cells.Item(7,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(7,4).Value = 44126.61624953821
$ws1.Cells.Item(7,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(7,7).Value = "No, individual is not listed"
$ws1.Cells.Item(8,3).Value = 43760.61627357056
$ws1.Cells.Item(8,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(8,4).Value = 44126.61627357056
$ws1.Cells.Item(8,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(8,7).Value = "No, individual is not listed"
$ws1.Cells.Item(9,3).Value = 43760.6162977645
$ws1.Cells.Item(9,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(9,4).Value = 44126.6162977645
$ws1.Cells.Item(9,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(9,7).Value = "No, individual is not listed"
$ws1.Cells.Item(10,3).Value = 43760.61632270875
$ws1.Cells.Item(10,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(10,4).Value = 44126.61632270875
$ws1.Cells.Item(10,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(10,7).Value = "No, individual is not listed"
$ws1.Cells.Item(11,3).Value = 43760.61634633726
$ws1.Cells.Item(11,3).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(11,4).Value = 44126.61634633726
$ws1.Cells.Item(11,4).NumberFormat = "yyyy\-mm\-dd\ h:mm:ss"
$ws1.Cells.Item(11,7).Value = "No, individual is not listed"

# Row 2 H column flips from "Yes" to "No" once the individual is re-checked
$ws1.Cells.Item(2,8).Value = "No, individual is not listed"

# --- Add sheet2 (rejected authors) right after Sheet1 ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "sheet2"
$ws2.Outline.SummaryRow = 1
$ws2.Outline.SummaryColumn = 1
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

$ws2.Cells.Item(1,1).Value = "hello world "
$ws2.Cells.Item(1,2).Value = "hello world"
$ws2.Cells.Item(1,3).Value = "hello world "
$ws2.Cells.Item(2,1).Value = "hello world "
$ws2.Cells.Item(3,1).Value = "hello world "
$ws2.Cells.Item(4,1).Value = "hello world "
$ws2.Cells.Item(5,1).Value = "hello world "
$ws2.Cells.Item(6,1).Value = "hello world "
$ws2.Cells.Item(7,1).Value = "hello world "
$ws2.Cells.Item(8,1).Value = "hello world "
$ws2.Cells.Item(9,1).Value = "hello world "
$ws2.Cells.Item(10,1).Value = "hello world "
$ws2.Cells.Item(11,1).Value = "hello world "
$ws2.Cells.Item(12,1).Value = "hello world "
$ws2.Cells.Item(13,1).Value = "hello world "
$ws2.Cells.Item(14,1).Value = "hello world "
$ws2.Cells.Item(15,1).Value = "hello world "
$ws2.Cells.Item(16,1).Value = "Hiiii"
$ws2.Cells.Item(17,1).Value = "Hiiii"
$ws2.Cells.Item(18,1).Value = "Moore"
$ws2.Cells.Item(19,1).Value = "Hello World "
$ws2.Cells.Item(20,1).Value = "Achiron"
$ws2.Cells.Item(21,1).Value = "Afsar"
$ws2.Cells.Item(22,1).Value = "Akgun"
$ws2.Cells.Item(23,1).Value = "Alroughani"
$ws2.Cells.Item(24,1).Value = "Bass"
$ws2.Cells.Item(25,1).Value = "Berkovich"
$ws2.Cells.Item(26,1).Value = "Broadley"
$ws2.Cells.Item(27,1).Value = "Celius"

# --- Add sheet21 and sheet22 (empty placeholder sheets) ---
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "sheet21"
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 0.75 * 72
$ws3.PageSetup.RightMargin = 0.75 * 72
$ws3.PageSetup.TopMargin = 1 * 72
$ws3.PageSetup.BottomMargin = 1 * 72
$ws3.PageSetup.HeaderMargin = 0.5 * 72
$ws3.PageSetup.FooterMargin = 0.5 * 72

$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "sheet22"
$ws4.Outline.SummaryRow = 1
$ws4.Outline.SummaryColumn = 1
$ws4.PageSetup.LeftMargin = 0.75 * 72
$ws4.PageSetup.RightMargin = 0.75 * 72
$ws4.PageSetup.TopMargin = 1 * 72
$ws4.PageSetup.BottomMargin = 1 * 72
$ws4.PageSetup.HeaderMargin = 0.5 * 72
$ws4.PageSetup.FooterMargin = 0.5 * 72

# --- Selections on each sheet ---
$ws1.Range("C15").Select() | Out-Null
$ws3.Range("A1").Select() | Out-Null
$ws4.Range("A1").Select() | Out-Null
$ws2.Range("E12").Select() | Out-Null

# --- Make sheet2 the active tab (matches activeTab=1 / tabSelected on sheet2) ---
$ws2.Activate() | Out-Null
